$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results for Case_3_29 (380 kV case)
$ws.Range("B2").Value2 = 13.88541484504435
$ws.Range("C2").Value2 = 5.177136891415603
$ws.Range("D2").Value2 = 8.833855649848598
$ws.Range("F2").Value2 = 44.30655457852706
$ws.Range("G2").Value2 = 3.739416386551787
$ws.Range("I2").Value2 = 36.47434230969247
$ws.Range("K2").Value2 = 12.10747298028357
$ws.Range("L2").Value2 = 11.03415644067324
$ws.Range("M2").Value2 = 15.63597724056909
$ws.Range("B3").Value2 = 13.80354887461018
$ws.Range("C3").Value2 = 4.968046568627969
$ws.Range("D3").Value2 = 8.822948447706572
$ws.Range("F3").Value2 = 43.96480458125354
$ws.Range("G3").Value2 = 3.742477462450222
$ws.Range("I3").Value2 = 36.30647604017987
$ws.Range("K3").Value2 = 12.04614633208203
$ws.Range("L3").Value2 = 11.03902050425428
$ws.Range("M3").Value2 = 15.64438632396853
$ws.Range("B4").Value2 = 13.75782338809849
$ws.Range("C4").Value2 = 4.833807079950862
$ws.Range("D4").Value2 = 8.815973170676974
$ws.Range("F4").Value2 = 43.75991431535026
$ws.Range("G4").Value2 = 3.744454971577026
$ws.Range("I4").Value2 = 36.20685230721862
$ws.Range("K4").Value2 = 12.01227811495379
$ws.Range("L4").Value2 = 11.04359143774879
$ws.Range("M4").Value2 = 15.65269541525561
$ws.Range("B5").Value2 = 13.74035141856255
$ws.Range("C5").Value2 = 4.777666547857674
$ws.Range("D5").Value2 = 8.813060433295707
$ws.Range("F5").Value2 = 43.67772037332951
$ws.Range("G5").Value2 = 3.745285551951484
$ws.Range("I5").Value2 = 36.16713944867798
$ws.Range("K5").Value2 = 11.99944215502729
$ws.Range("L5").Value2 = 11.04585290051465
$ws.Range("M5").Value2 = 15.65687307592258
$ws.Range("B6").Value2 = 13.73752090213602
$ws.Range("C6").Value2 = 4.768258827899938
$ws.Range("D6").Value2 = 8.812572523446702
$ws.Range("F6").Value2 = 43.66415229194309
$ws.Range("G6").Value2 = 3.745424965351138
$ws.Range("I6").Value2 = 36.16059908766457
$ws.Range("K6").Value2 = 11.99736943603819
$ws.Range("L6").Value2 = 11.04625250810413
$ws.Range("M6").Value2 = 15.65761459741077
$ws.Range("B7").Value2 = 13.75758302743721
$ws.Range("C7").Value2 = 4.833055715434798
$ws.Range("D7").Value2 = 8.815934173349605
$ws.Range("F7").Value2 = 43.75880048026967
$ws.Range("G7").Value2 = 3.744466072831117
$ws.Range("I7").Value2 = 36.20631312125705
$ws.Range("K7").Value2 = 12.01210107841762
$ws.Range("L7").Value2 = 11.04362032161398
$ws.Range("M7").Value2 = 15.65274855080376
$ws.Range("B8").Value2 = 13.85625756931396
$ws.Range("C8").Value2 = 5.106278207922705
$ws.Range("D8").Value2 = 8.830152410325571
$ws.Range("F8").Value2 = 44.18772757567174
$ws.Range("G8").Value2 = 3.74045155914135
$ws.Range("I8").Value2 = 36.41575768428031
$ws.Range("K8").Value2 = 12.08555019053257
$ws.Range("L8").Value2 = 11.03550497848291
$ws.Range("M8").Value2 = 15.63822419897297
$ws.Range("B9").Value2 = 14.08481689159845
$ws.Range("C9").Value2 = 5.594265066301577
$ws.Range("D9").Value2 = 8.85584048633187
$ws.Range("F9").Value2 = 45.0654905555968
$ws.Range("G9").Value2 = 3.733352676677505
$ws.Range("I9").Value2 = 36.8529982759777
$ws.Range("K9").Value2 = 12.25896320958884
$ws.Range("L9").Value2 = 11.03214331617741
$ws.Range("M9").Value2 = 15.63466681889311
$ws.Range("B10").Value2 = 14.27268343601058
$ws.Range("C10").Value2 = 5.922184962427582
$ws.Range("D10").Value2 = 8.873395296671536
$ws.Range("F10").Value2 = 45.72910486230094
$ws.Range("G10").Value2 = 3.728603131154698
$ws.Range("I10").Value2 = 37.18923819076663
$ws.Range("K10").Value2 = 12.40328807448879
$ws.Range("L10").Value2 = 11.03729639017278
$ws.Range("M10").Value2 = 15.64718323627674
$ws.Range("B11").Value2 = 14.36214089985726
$ws.Range("C11").Value2 = 6.064486679814562
$ws.Range("D11").Value2 = 8.881099263189654
$ws.Range("F11").Value2 = 46.03426382672212
$ws.Range("G11").Value2 = 3.726542437586641
$ws.Range("I11").Value2 = 37.34519065000958
$ws.Range("K11").Value2 = 12.47238007679264
$ws.Range("L11").Value2 = 11.0412874692555
$ws.Range("M11").Value2 = 15.65614228359717
$ws.Range("B12").Value2 = 14.39656034159565
$ws.Range("C12").Value2 = 6.117369318971424
$ws.Range("D12").Value2 = 8.883976333551438
$ws.Range("F12").Value2 = 46.15022155391897
$ws.Range("G12").Value2 = 3.725776380904286
$ws.Range("I12").Value2 = 37.40465217610114
$ws.Range("K12").Value2 = 12.49901604456159
$ws.Range("L12").Value2 = 11.04303467026
$ws.Range("M12").Value2 = 15.66000206450005
$ws.Range("B13").Value2 = 14.38912381286174
$ws.Range("C13").Value2 = 6.106024971203334
$ws.Range("D13").Value2 = 8.883358491064888
$ws.Range("F13").Value2 = 46.12523129482398
$ws.Range("G13").Value2 = 3.725940730966969
$ws.Range("I13").Value2 = 37.39182843443288
$ws.Range("K13").Value2 = 12.49325886638459
$ws.Range("L13").Value2 = 11.04264790558506
$ws.Range("M13").Value2 = 15.65915005461239
$ws.Range("B14").Value2 = 14.36496188875662
$ws.Range("C14").Value2 = 6.068857568162746
$ws.Range("D14").Value2 = 8.881336766737137
$ws.Range("F14").Value2 = 46.04379611593919
$ws.Range("G14").Value2 = 3.726479127852905
$ws.Range("I14").Value2 = 37.3500746003434
$ws.Range("K14").Value2 = 12.47456209736385
$ws.Range("L14").Value2 = 11.04142649134038
$ws.Range("M14").Value2 = 15.65645048173418
$ws.Range("B15").Value2 = 14.35023190486968
$ws.Range("C15").Value2 = 6.045960293206758
$ws.Range("D15").Value2 = 8.88009316778207
$ws.Range("F15").Value2 = 45.99396478653292
$ws.Range("G15").Value2 = 3.726810769165779
$ws.Range("I15").Value2 = 37.32455125226046
$ws.Range("K15").Value2 = 12.46317062697374
$ws.Range("L15").Value2 = 11.04070902608637
$ws.Range("M15").Value2 = 15.65485768082639
$ws.Range("B16").Value2 = 14.26691498459742
$ws.Range("C16").Value2 = 5.912745528258938
$ws.Range("D16").Value2 = 8.872886176759872
$ws.Range("F16").Value2 = 45.70922193389225
$ws.Range("G16").Value2 = 3.728739805515866
$ws.Range("I16").Value2 = 37.17910448373631
$ws.Range("K16").Value2 = 12.39884006707783
$ws.Range("L16").Value2 = 11.03706862470971
$ws.Range("M16").Value2 = 15.64666325320147
$ws.Range("B17").Value2 = 14.21680468947868
$ws.Range("C17").Value2 = 5.829250593848502
$ws.Range("D17").Value2 = 8.86839293660697
$ws.Range("F17").Value2 = 45.53533082021609
$ws.Range("G17").Value2 = 3.729948734715691
$ws.Range("I17").Value2 = 37.09062722166485
$ws.Range("K17").Value2 = 12.36024051239285
$ws.Range("L17").Value2 = 11.03525652702186
$ws.Range("M17").Value2 = 15.64247098062624
$ws.Range("B18").Value2 = 14.18836061859751
$ws.Range("C18").Value2 = 5.780580575055781
$ws.Range("D18").Value2 = 8.865782003098218
$ws.Range("F18").Value2 = 45.43562602549078
$ws.Range("G18").Value2 = 3.730653486614445
$ws.Range("I18").Value2 = 37.04002091506272
$ws.Range("K18").Value2 = 12.33836404633424
$ws.Range("L18").Value2 = 11.0343693371328
$ws.Range("M18").Value2 = 15.64036720126229
$ws.Range("B19").Value2 = 14.17879579414981
$ws.Range("C19").Value2 = 5.763991333965127
$ws.Range("D19").Value2 = 8.864893419248943
$ws.Range("F19").Value2 = 45.40192355566759
$ws.Range("G19").Value2 = 3.730893721802406
$ws.Range("I19").Value2 = 37.02293592479858
$ws.Range("K19").Value2 = 12.3310135213689
$ws.Range("L19").Value2 = 11.03409561243521
$ws.Range("M19").Value2 = 15.63970777823278
$ws.Range("B20").Value2 = 14.2221001092392
$ws.Range("C20").Value2 = 5.838205748784167
$ws.Range("D20").Value2 = 8.868873994596255
$ws.Range("F20").Value2 = 45.55380996502045
$ws.Range("G20").Value2 = 3.729819069066306
$ws.Range("I20").Value2 = 37.10001661747788
$ws.Range("K20").Value2 = 12.36431602366656
$ws.Range("L20").Value2 = 11.03543338423399
$ws.Range("M20").Value2 = 15.64288544566175
$ws.Range("B21").Value2 = 14.37204432245923
$ws.Range("C21").Value2 = 6.07980189518286
$ws.Range("D21").Value2 = 8.881931686554925
$ws.Range("F21").Value2 = 46.06770529451062
$ws.Range("G21").Value2 = 3.726320600644605
$ws.Range("I21").Value2 = 37.36232790558853
$ws.Range("K21").Value2 = 12.48004115316237
$ws.Range("L21").Value2 = 11.04177885720414
$ws.Range("M21").Value2 = 15.65723075294621
$ws.Range("B22").Value2 = 14.47319647460252
$ws.Range("C22").Value2 = 6.231841737244979
$ws.Range("D22").Value2 = 8.890231130613198
$ws.Range("F22").Value2 = 46.40586776589978
$ws.Range("G22").Value2 = 3.724117365263377
$ws.Range("I22").Value2 = 37.53611619144147
$ws.Range("K22").Value2 = 12.55841643595024
$ws.Range("L22").Value2 = 11.04730026996979
$ws.Range("M22").Value2 = 15.66932801517274
$ws.Range("B23").Value2 = 14.41893142463049
$ws.Range("C23").Value2 = 6.151235714592501
$ws.Range("D23").Value2 = 8.885822931297426
$ws.Range("F23").Value2 = 46.22519692384521
$ws.Range("G23").Value2 = 3.725285686168244
$ws.Range("I23").Value2 = 37.44315529334899
$ws.Range("K23").Value2 = 12.51634267656352
$ws.Range("L23").Value2 = 11.04422798868453
$ws.Range("M23").Value2 = 15.66262330217727
$ws.Range("B24").Value2 = 14.21970491069654
$ws.Range("C24").Value2 = 5.834159201687124
$ws.Range("D24").Value2 = 8.868656594604046
$ws.Range("F24").Value2 = 45.5454547049178
$ws.Range("G24").Value2 = 3.72987766061566
$ws.Range("I24").Value2 = 37.09577085841634
$ws.Range("K24").Value2 = 12.36247250238418
$ws.Range("L24").Value2 = 11.03535294543642
$ws.Range("M24").Value2 = 15.64269711140643
$ws.Range("B25").Value2 = 14.0193729393401
$ws.Range("C25").Value2 = 5.467535764169174
$ws.Range("D25").Value2 = 8.849124047095046
$ws.Range("F25").Value2 = 44.82447616590586
$ws.Range("G25").Value2 = 3.735190875811095
$ws.Range("I25").Value2 = 36.73199043965541
$ws.Range("K25").Value2 = 12.20900214093454
$ws.Range("L25").Value2 = 11.03171149984195
$ws.Range("M25").Value2 = 15.63296621154018
